# Apply updated cryptocurrency price/volume data as described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.641.82"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "1.593.33"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "'210.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.20%  "
$ws.Range("D6").Value = "'0.516"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.06%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").Value = "'0.0616"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("E9").Value = "  -1.69%  "
$ws.Range("D10").Value = "'19.44"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.18%  "
$ws.Range("D11").Value = "'0.0839"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.46%  "
$ws.Range("D12").Value = "1.815.87"
$ws.Range("E12").Value = "  +0.15%  "
$ws.Range("D13").Value = "1.583.02"
$ws.Range("E13").Value = "  -0.23%  "
$ws.Range("D14").Value = "'4.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.44%  "
$ws.Range("D15").Value = "'0.521"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.41%  "
$ws.Range("D16").Value = "'64.43"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.25%  "
$ws.Range("D17").Value = "26.605.04"
$ws.Range("E17").Value = "  -0.17%  "
$ws.Range("D18").Value = "0.0₃0728"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("D20").Value = "'207.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.33%  "
$ws.Range("D21").Value = "'6.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.70%  "
$ws.Range("E22").Value = "  -0.31%  "
$ws.Range("E23").Value = "  -2.80%  "
$ws.Range("D24").Value = "'8.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.54%  "
$ws.Range("D25").Value = "'145.55"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.18%  "
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("E27").Value = "  -2.47%  "
$ws.Range("E28").Value = "  +0.37%  "
$ws.Range("D29").Value = "'15.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("E31").Value = "  +0.14%  "
$ws.Range("D32").Value = "'3.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.52%  "
$ws.Range("D33").Value = "'0.657"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.99%  "
$ws.Range("D34").Value = "'2.91"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.16%  "
$ws.Range("D35").Value = "1.277.97"
$ws.Range("E35").Value = "  -3.88%  "
$ws.Range("E36").Value = "  +1.78%  "
$ws.Range("D37").Value = "'1.49"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.01%  "
$ws.Range("E38").Value = "  -0.32%  "
$ws.Range("D39").Value = "'0.837"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.17%  "
$ws.Range("E40").Value = "  +0.20%  "
$ws.Range("D41").Value = "'5.42"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.46%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'0.787"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.32%  "
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").Value = "'2.19"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.79%  "
$ws.Range("D44").Value = "'63.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").Value = "'0.915"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +9.31%  "
$ws.Range("D46").Value = "1.728.31"
$ws.Range("E46").Value = "  +0.12%  "
$ws.Range("D47").Value = "'89.71"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.25%  "
$ws.Range("E48").Value = "  -0.44%  "
$ws.Range("E49").Value = "  -3.36%  "
$ws.Range("E50").Value = "  +3.00%  "
$ws.Range("E51").Value = "  -0.96%  "
